$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Shade the header row of "Table 1: Sensor information" (3 cells) and the
#    header row of "Table 2: Control of sensor readings" (3 cells) with the
#    light-blue accent1/33% theme fill used for the new table headers.
# ---------------------------------------------------------------------------
$fillColor = 16114881   # BGR-encoded 0xF5E4C1 -> renders as w:fill="C1E4F5"
$autoColor = -16777216  # wdColorAutomatic -> renders as w:color="auto"

$tablesToShade = @(1, 2)
foreach ($tIdx in $tablesToShade) {
    $tbl = $d.Tables.Item($tIdx)
    $headerRow = $tbl.Rows.Item(1)
    $colCount = $headerRow.Cells.Count
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $headerRow.Cells.Item($c)
        $cell.Shading.Texture = 0
        $cell.Shading.ForegroundPatternColor = $autoColor
        $cell.Shading.BackgroundPatternColor = $fillColor
    }
}

# ---------------------------------------------------------------------------
# 2) Split "The information from the sensors is sent to Node-Red ..." so the
#    word "is" becomes its own run (the sentence's text is unchanged).
#    NB: use $d.Content.Paragraphs (derived fresh from the whole-document
#    range) rather than the cached $d.Paragraphs, since touching
#    $d.Tables.Item(...) above can leave $d.Paragraphs scoped to that table.
# ---------------------------------------------------------------------------
$paraCount = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Content.Paragraphs.Item($i)
    $full = $p.Range.Text
    if ($full.StartsWith("The information from the sensors is sent to Node-Red")) {
        $pStart = $p.Range.Start
        $idx = $full.IndexOf("is sent to Node-Red")
        $isStart = $pStart + $idx
        $isEnd = $isStart + 2
        $isRange = $d.Range($isStart, $isEnd)
        $isRange.Font.Bold = 1
        $isRange.Font.Bold = 0
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Collapse the "[12]" bibliography entry (currently split across several
#    runs) into a single run with the same, unchanged text.
# ---------------------------------------------------------------------------
$citation = "[12] Fundamentals of Environmental Measurements. “Dissolved Oxygen.” fondriest.com, https://www.fondriest.com/environmental-measurements/parameters/water-quality/dissolved-oxygen/ (accessed Apr. 2, 2024)."

$lastPara = $d.Content.Paragraphs.Last
$lastRange = $lastPara.Range
$start = $lastRange.Start
$end = $lastRange.End
$bodyRange = $d.Range($start, $end - 1)
if ($bodyRange.Text -eq "[12] Fundamentals of Environmental Measurements. " + [char]0x201C + "Dissolved Oxygen." + [char]0x201D + " fondriest.com, https://www.fondriest.com/environmental-measurements/parameters/water-quality/dissolved-oxygen/ (accessed Apr. 2, 2024).") {
    # Force a full re-segmentation: first collapse to transient placeholder
    # text (runs merge into one when content genuinely changes), then write
    # back the real text so the paragraph ends up as a single clean run.
    $bodyRange.Text = "~"
    $lastPara2 = $d.Content.Paragraphs.Last
    $lastRange2 = $lastPara2.Range
    $start2 = $lastRange2.Start
    $end2 = $lastRange2.End
    $bodyRange2 = $d.Range($start2, $end2 - 1)
    $bodyRange2.Text = $citation
}

Write-Host "edits applied"
